# Auto-generated edit script applying the Masamune_Profits.xlsx diff
# Updates Leve profit-calculation columns (H-N) across the ALC, ARM, BSM, CRP, GSM, and LTW sheets
# to reflect refreshed market-board pricing data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 14493931
$ws.Range("I80").Value = 47619892
$ws.Range("J80").Value = 1322.25
$ws.Range("K80").Value = 142859676
$ws.Range("L80").Value = 3966.75
$ws.Range("M80").Value = -142858678
$ws.Range("N80").Value = -5962.75
$ws.Range("H83").Value = 14493931
$ws.Range("I83").Value = 47619892
$ws.Range("J83").Value = 1322.25
$ws.Range("K83").Value = 428579028
$ws.Range("L83").Value = 11900.25
$ws.Range("M83").Value = -428574036
$ws.Range("N83").Value = -21884.25
$ws.Range("H137").Value = 3082346.5
$ws.Range("I137").Value = 10990037
$ws.Range("J137").Value = 7133.4443
$ws.Range("K137").Value = 32970111
$ws.Range("L137").Value = 21400.3329
$ws.Range("M137").Value = -32967561
$ws.Range("N137").Value = -26500.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 6336
$ws.Range("I19").Value = 3004
$ws.Range("J19").Value = 13000
$ws.Range("K19").Value = 3004
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = -2775
$ws.Range("N19").Value = -13458
$ws.Range("H32").Value = 12313
$ws.Range("I32").Value = 10899.191
$ws.Range("J32").Value = 25602.8
$ws.Range("K32").Value = 10899.191
$ws.Range("L32").Value = 25602.8
$ws.Range("M32").Value = -10612.191
$ws.Range("N32").Value = -26176.8
$ws.Range("H36").Value = 1150
$ws.Range("I36").Value = 1150
$ws.Range("K36").Value = 1150
$ws.Range("M36").Value = -804
$ws.Range("H64").Value = 29560
$ws.Range("J64").Value = 29560
$ws.Range("L64").Value = 29560
$ws.Range("N64").Value = -30056
$ws.Range("H67").Value = 29560
$ws.Range("J67").Value = 29560
$ws.Range("L67").Value = 29560
$ws.Range("N67").Value = -31276
$ws.Range("H74").Value = 2206.9375
$ws.Range("I74").Value = 2002.3334
$ws.Range("K74").Value = 2002.3334
$ws.Range("M74").Value = -1128.3334
$ws.Range("H77").Value = 2206.9375
$ws.Range("I77").Value = 2002.3334
$ws.Range("K77").Value = 10011.667
$ws.Range("M77").Value = -5643.666999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 417.66666
$ws.Range("I5").Value = 124
$ws.Range("J5").Value = 1005
$ws.Range("K5").Value = 124
$ws.Range("L5").Value = 1005
$ws.Range("M5").Value = -11
$ws.Range("N5").Value = -1231
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
$ws.Range("H141").Value = 37500
$ws.Range("J141").Value = 37500
$ws.Range("L141").Value = 37500
$ws.Range("N141").Value = -47860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 595.9091
$ws.Range("I19").Value = 283.8889
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 283.8889
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = -113.8889
$ws.Range("N19").Value = -2340
$ws.Range("H24").Value = 595.9091
$ws.Range("I24").Value = 283.8889
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 283.8889
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = -113.8889
$ws.Range("N24").Value = -2340
$ws.Range("H31").Value = 7942739
$ws.Range("I31").Value = 2506.85
$ws.Range("J31").Value = 15161132
$ws.Range("K31").Value = 2506.85
$ws.Range("L31").Value = 15161132
$ws.Range("M31").Value = -2211.85
$ws.Range("N31").Value = -15161722
$ws.Range("H34").Value = 7942739
$ws.Range("I34").Value = 2506.85
$ws.Range("J34").Value = 15161132
$ws.Range("K34").Value = 2506.85
$ws.Range("L34").Value = 15161132
$ws.Range("M34").Value = -2304.85
$ws.Range("N34").Value = -15161536
$ws.Range("H132").Value = 826251.3
$ws.Range("I132").Value = 1840.3334
$ws.Range("K132").Value = 5521.0002
$ws.Range("M132").Value = -2991.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2409657.5
$ws.Range("I2").Value = 3012059.2
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 3012059.2
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = -3011946.2
$ws.Range("N2").Value = -276
$ws.Range("H18").Value = 11250
$ws.Range("I18").Value = 10000
$ws.Range("J18").Value = 15000
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = -9707
$ws.Range("N18").Value = -15586
$ws.Range("H31").Value = 11434.333
$ws.Range("I31").Value = 1654
$ws.Range("J31").Value = 30995
$ws.Range("K31").Value = 1654
$ws.Range("L31").Value = 30995
$ws.Range("M31").Value = -1362
$ws.Range("N31").Value = -31579
$ws.Range("H37").Value = 11434.333
$ws.Range("I37").Value = 1654
$ws.Range("J37").Value = 30995
$ws.Range("K37").Value = 1654
$ws.Range("L37").Value = 30995
$ws.Range("M37").Value = -1377
$ws.Range("N37").Value = -31549
$ws.Range("H43").Value = 9000
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = -2849
$ws.Range("N43").Value = -15302
$ws.Range("H57").Value = 13469.091
$ws.Range("I57").Value = 6318.1816
$ws.Range("J57").Value = 20620
$ws.Range("K57").Value = 6318.1816
$ws.Range("L57").Value = 20620
$ws.Range("M57").Value = -5498.1816
$ws.Range("N57").Value = -22260
$ws.Range("H70").Value = 5096.8486
$ws.Range("I70").Value = 5191.9165
$ws.Range("J70").Value = 4843.3335
$ws.Range("K70").Value = 5191.9165
$ws.Range("L70").Value = 4843.3335
$ws.Range("M70").Value = -4921.9165
$ws.Range("N70").Value = -5383.3335
$ws.Range("H73").Value = 5096.8486
$ws.Range("I73").Value = 5191.9165
$ws.Range("J73").Value = 4843.3335
$ws.Range("K73").Value = 5191.9165
$ws.Range("L73").Value = 4843.3335
$ws.Range("M73").Value = -4255.9165
$ws.Range("N73").Value = -6715.3335
$ws.Range("H80").Value = 6410.8423
$ws.Range("I80").Value = 4733.3335
$ws.Range("J80").Value = 9286.571
$ws.Range("K80").Value = 4733.3335
$ws.Range("L80").Value = 9286.571
$ws.Range("M80").Value = -3735.3335
$ws.Range("N80").Value = -11282.571
$ws.Range("H83").Value = 6410.8423
$ws.Range("I83").Value = 4733.3335
$ws.Range("J83").Value = 9286.571
$ws.Range("K83").Value = 23666.6675
$ws.Range("L83").Value = 46432.855
$ws.Range("M83").Value = -18674.6675
$ws.Range("N83").Value = -56416.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2833.5
$ws.Range("I9").Value = 238.28572
$ws.Range("J9").Value = 21000
$ws.Range("K9").Value = 238.28572
$ws.Range("L9").Value = 21000
$ws.Range("M9").Value = -14.28572
$ws.Range("N9").Value = -21448
$ws.Range("H82").Value = 27778612
$ws.Range("I82").Value = 1502
$ws.Range("J82").Value = 41667170
$ws.Range("K82").Value = 1502
$ws.Range("L82").Value = 41667170
$ws.Range("M82").Value = -1141
$ws.Range("N82").Value = -41667892
$ws.Range("H85").Value = 27778612
$ws.Range("I85").Value = 1502
$ws.Range("J85").Value = 41667170
$ws.Range("K85").Value = 1502
$ws.Range("L85").Value = 41667170
$ws.Range("M85").Value = -254
$ws.Range("N85").Value = -41669666
